# Update Name of Algo - correct a handful of imputed values in the result data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value  = -12.7646
$ws.Range("A8").Value  = -21.12650000000001
$ws.Range("A10").Value = -20.54039999999998
$ws.Range("A12").Value = -22.41430000000004
$ws.Range("B13").Value = 6.046000000000001
$ws.Range("A18").Value = -22.32020000000003
$ws.Range("C20").Value = -13.5421
